$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-3: update values to "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# Insert 10 new rows right after row 3 (before current row 4),
# then populate them with the new benchmark values.
$beforeRow = $t.Rows.Item(4)
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Add($beforeRow) | Out-Null

$newValues = @("104", "0.00003", "0.00016", "0.00005", "0.00002", "0.00004", "0.00005", "0.00009", "0.00470", "100.0")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $t.Cell(4 + $i, 1).Range.Text = $newValues[$i]
}

# The last three rows (originally 34-36, now shifted by +10 to 44-46) had
# multi-run tab-separated content; collapse each back down to a single value.
$t.Cell(44, 1).Range.Text = "99.97"
$t.Cell(45, 1).Range.Text = "0"
$t.Cell(46, 1).Range.Text = "15"
